$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New commit row: "Goblin added and creature functions update", 2 hours
$ws.Range("C15").Value = "Goblin added and creature functions update"
$ws.Range("C15").Style = $wb.Styles.Item("20% - Énfasis5")
$ws.Range("G15").Value = 2

# Selection / view ends up on C15 (matches author last click before saving)
$ws.Range("C15").Select() | Out-Null

# Page setup: A4, portrait
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "edit applied"
